$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Include from hp.owl")

# Remove the "Pediatric (<= 15 years)" concept row (HP:0410280); everything below shifts up one row.
$ws.Rows.Item(4).Delete()

# Temporarily remove the blank separator row and the "System URI" row so they get
# re-appended (in their original relative order) after the new concept rows below.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(9).Delete()

# Rename "Juvenile " to the fully-qualified age range label.
$ws.Range("B7").Value = "Juvenile (>= 5 years and < 16 years)"

# Replace the old single "Adult" bucket with the new "Young Adult" bucket.
$ws.Range("A8").Value = "HP:0011462"
$ws.Range("B8").Value = "Young Adult (>= 16 years and < 40 years)"

# Add new concept rows for "Middle Age" and "Senior".
$ws.Range("A9").Value = "HP:0003596"
$ws.Range("B9").Value = "Middle Age (>= 40 years and < 60 years)"
$ws.Range("A10").Value = "HP:0003584"
$ws.Range("B10").Value = "Senior (>= 60 years)"

# Restore the blank separator row and the "System URI" row at the bottom of the table.
$ws.Range("A12").Value = "System URI"
$ws.Range("B12").Value = "http://purl.obolibrary.org/obo/hp.owl"
